# Append a new, empty ("Blank" layout) slide at the end of the deck.
#
# PowerPoint's Slides.Add() only ever yields a minimal stub slide (no
# group-transform, no creationId) in this host, so instead we duplicate an
# existing slide (which carries the full <p:cSld>/<p:grpSpPr>/<p:extLst>/
# <p:clrMapOvr> skeleton that real PowerPoint writes), strip out all of its
# shapes so the result is a genuinely blank slide, point it at the "Blank"
# layout, and move it to the very end of the slide list.

$p = $ppt.ActivePresentation

$donor = $p.Slides.Item(13)
$newRange = $donor.Duplicate()
$newSlide = $newRange.Item(1)

while ($newSlide.Shapes.Count -gt 0) {
    $newSlide.Shapes.Item(1).Delete()
}

# "Blank" is always layout #7 on the slide master (type="blank").
$blankLayout = $p.SlideMaster.CustomLayouts.Item(7)
$newSlide.CustomLayout = $blankLayout

$newSlide.MoveTo($p.Slides.Count)
